$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths ---
# Target stored widths are 47 (col B) and 30 (col D). This engine's
# ColumnWidth setter re-adds a constant ~5/6 character padding when it
# serializes back to the <col width="..."> attribute, so back that out here
# so the saved file ends up with exactly width="47" / width="30".
$ws.Columns.Item(2).ColumnWidth = 47 - (5/6)
$ws.Columns.Item(4).ColumnWidth = 30 - (5/6)

# --- Rows 2 & 3: unchanged except the refreshed "captured at" timestamp ---
$ws.Range("A2").Value = '2025-12-28 01:43:02'
$ws.Range("A3").Value = '2025-12-28 01:43:02'

# --- A new listing lands at row 4 (pushing the former rows 4-5 down to
# --- rows 5-6), and two more new listings are appended as rows 7-8.
# --- Rewrite rows 4-8 in full with the post-append content. ---
$ws.Range("A4").Value = '2025-12-28 01:43:02'
$ws.Range("B4").Value = '【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5217096'
$ws.Range("G4").Value = 243
$ws.Range("H4").Value = '🔥API ◆ツール'

$ws.Range("A5").Value = '2025-12-28 01:43:02'
$ws.Range("B5").Value = '【募集】Web予約フロー自動化ツールの設計・開発をお任せします'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5462249'
$ws.Range("G5").Value = 213
$ws.Range("H5").Value = '◆ツール,開発'

$ws.Range("A6").Value = '2025-12-28 01:43:02'
$ws.Range("B6").Value = '【急募】ReactでLine風会話履歴表示コンポーネント作成依頼'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5462198'
$ws.Range("G6").Value = 128
$ws.Range("H6").Value = '🔥React'

$ws.Range("A7").Value = '2025-12-28 01:43:02'
$ws.Range("B7").Value = 'GoogleCloudを利用したアジャイル開発共通基盤のSREエンジニアの募集'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5457458'
$ws.Range("G7").Value = 75
$ws.Range("H7").Value = '◆開発'

$ws.Range("A8").Value = '2025-12-28 01:43:02'
$ws.Range("B8").Value = 'FXレイテンシーアービトラージの検証(環境設計・比較評価・PoC)'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5462397'
$ws.Range("G8").Value = 25
# (row 8 has no skill-summary column H, same as the diff)

# --- Hyperlinks ---
# This runtime's Range.Hyperlinks.Delete() clears every hyperlink on the
# sheet (it isn't scoped to the calling range), so drop them all first and
# re-add F2:F8 in the final top-to-bottom order. That reproduces rId1..rId7
# in exactly the order the diff's relationship part lists them in.
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5462048')
$ws.Range("F2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5461891')
$ws.Range("F3").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5217096')
$ws.Range("F4").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5462249')
$ws.Range("F5").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5462198')
$ws.Range("F6").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5457458')
$ws.Range("F7").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5462397')
$ws.Range("F8").Style = "Hyperlink"

Write-Output "applied"
